$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1K_ADAMTS19_01", 0, 0, 201),
    @("1K_ADAMTS19_03", 0, 2, 286),
    @("1K_ADAMTS19_04", 0, 0, 154),
    @("1K_ADAMTS19_10", 0, 0, 21),
    @("1K_ADAMTS19_10", 1, 0, 16),
    @("1K_ADAMTS19_10", 2, 0, 22),
    @("1K_ADAMTS19_10", 3, 0, 17),
    @("1K_ADAMTS19_10", 4, 0, 21),
    @("1K_ADAMTS19_10", 5, 0, 36),
    @("1K_ADAMTS19_10", 6, 0, 27),
    @("1K_ADAMTS19_10", 7, 0, 30),
    @("1K_ADAMTS19_10", 8, 0, 32),
    @("1K_ADAMTS19_10", 9, 0, 15),
    @("1K_ADAMTS19_10", 10, 0, 15),
    @("1K_ADAMTS19_10", 11, 0, 25),
    @("1K_SCR_01", 0, 0, 14),
    @("1K_SCR_01", 1, 0, 19),
    @("1K_SCR_03", 0, 1, 15),
    @("1K_SCR_03", 1, 1, 60),
    @("1K_SCR_03", 2, 1, 13),
    @("1K_SCR_03", 3, 1, 12),
    @("1K_SCR_03", 4, 1, 29),
    @("1K_SCR_03", 5, 1, 175),
    @("1K_SCR_03", 6, 1, 553),
    @("1K_SCR_03", 7, 1, 12),
    @("1K_TAX1BP1_01", 0, 0, 151),
    @("1K_TAX1BP1_01", 1, 0, 16),
    @("1K_TAX1BP1_01", 2, 0, 14),
    @("1K_TAX1BP1_01", 3, 0, 18),
    @("1K_TAX1BP1_02", 0, 0, 11),
    @("1K_TAX1BP1_02", 1, 0, 22),
    @("1K_TAX1BP1_02", 2, 0, 83),
    @("1K_TAX1BP1_02", 3, 0, 16),
    @("1K_TAX1BP1_02", 4, 0, 125),
    @("1K_TAX1BP1_02", 5, 0, 307),
    @("1K_TAX1BP1_02", 6, 0, 89),
    @("1K_TAX1BP1_02", 7, 0, 241),
    @("1K_TAX1BP1_02", 8, 0, 41),
    @("1K_TAX1BP1_02", 9, 0, 24),
    @("1K_TAX1BP1_02", 10, 0, 19),
    @("1K_TAX1BP1_02", 11, 0, 34),
    @("1K_TAX1BP1_02", 12, 0, 15),
    @("1K_TAX1BP1_02", 13, 0, 49),
    @("1K_TAX1BP1_02", 14, 0, 36),
    @("1K_TAX1BP1_02", 15, 0, 32),
    @("1K_TAX1BP1_02", 16, 0, 62),
    @("1K_TAX1BP1_02", 17, 0, 20),
    @("1K_TAX1BP1_02", 18, 0, 19),
    @("1K_TAX1BP1_02", 19, 0, 144),
    @("1K_TAX1BP1_02", 20, 0, 36),
    @("1K_TAX1BP1_02", 21, 0, 29),
    @("1K_TAX1BP1_02", 22, 0, 44),
    @("1K_TAX1BP1_02", 23, 0, 12),
    @("1K_TAX1BP1_02", 24, 0, 20),
    @("1K_TAX1BP1_07", 0, 0, 994),
    @("1K_TAX1BP1_08", 0, 1, 189),
    @("1K_TAX1BP1_08", 1, 1, 34),
    @("1K_TAX1BP1_09", 0, 0, 257),
    @("1K_TAX1BP1_10", 0, 7, 580),
    @("WT_ADAMTS19_01", 0, 0, 16),
    @("WT_ADAMTS19_01", 1, 0, 559),
    @("WT_ADAMTS19_02", 0, 1, 24),
    @("WT_ADAMTS19_02", 1, 1, 156),
    @("WT_ADAMTS19_02", 2, 1, 15),
    @("WT_ADAMTS19_02", 3, 1, 14),
    @("WT_ADAMTS19_02", 4, 1, 18),
    @("WT_ADAMTS19_07", 0, 0, 349),
    @("WT_ADAMTS19_07", 1, 0, 18),
    @("WT_ADAMTS19_07", 2, 0, 34),
    @("WT_ADAMTS19_07", 3, 0, 35),
    @("WT_ADAMTS19_10", 0, 1, 1472),
    @("WT_SCR_01", 0, 1, 519),
    @("WT_SCR_01", 1, 1, 33),
    @("WT_SCR_10", 0, 1, 541),
    @("WT_TAX1BP1_01", 0, 3, 301),
    @("WT_TAX1BP1_01", 1, 3, 1231),
    @("WT_TAX1BP1_01", 2, 3, 40),
    @("WT_TAX1BP1_02", 0, 1, 43),
    @("WT_TAX1BP1_03", 0, 0, 153),
    @("WT_TAX1BP1_03", 0, 2, 356),
    @("WT_TAX1BP1_04", 0, 0, 50),
    @("WT_TAX1BP1_04", 1, 0, 42),
    @("WT_TAX1BP1_04", 2, 0, 2607),
    @("WT_TAX1BP1_04", 0, 1, 1015),
    @("WT_TAX1BP1_04", 1, 1, 92),
    @("WT_TAX1BP1_06", 0, 2, 1150),
    @("WT_TAX1BP1_06", 1, 2, 264),
    @("WT_TAX1BP1_07", 0, 0, 286),
    @("WT_TAX1BP1_08", 0, 3, 373),
    @("WT_TAX1BP1_09", 0, 0, 308)
)

$rowCount = $data.Count
$colCount = 4

$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1

$range = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$range.Value2 = $values
